# The author re-uploaded the workbook after lower-casing the header row
# labels on the only worksheet ("Feuil1"):
#   A1 "Nom pilote" -> "nom pilote"
#   B1 "Abandon"    -> "abandon"
#   C1 "Grand prix" -> "grand prix"
#   M1 "écurie"     -> "ecurie"   (accent dropped too)
# All other header cells (D1:L1) were already lower-case and are untouched.
# No data cell values changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "nom pilote"
$ws.Range("B1").Value = "abandon"
$ws.Range("C1").Value = "grand prix"
$ws.Range("M1").Value = "ecurie"

# Reset the view back to the top-left / first cell (the saved file had
# scrolled to C1 with R9 selected).
$ws.Range("A1").Select()
